$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header label change ---
$ws.Range("B1").Value = "State 1"

# --- Row 2: updated energy value ---
$ws.Range("B2").Value = -149.562132347448

# --- Row 5: section header label swap (was "Psi4 (1x)", now "SF-XCIS") ---
$ws.Range("A5").Value = "SF-XCIS"

# --- Row 6: updated energy value; #dets column cleared ---
$ws.Range("B6").Value = -149.604321051649
$ws.Range("C6").ClearContents() | Out-Null

# --- Row 8: section header label swap (was "SF-XCIS", now "Psi4 (1x)") ---
$ws.Range("A8").Value = "Psi4 (1x)"

# --- Row 9: updated energy value; #dets column filled in ---
$ws.Range("B9").Value = -149.562132340033
$ws.Range("C9").Value = 108

# --- Row 10: new N2 data row under the "Psi4 (1x)" section ---
$ws.Range("A10").Value = "N2"
$ws.Range("B10").Value = -108.780110348207

# --- Row 12: new section header "Psi4 (S)" (bold, like the other section headers) ---
$ws.Range("A12").Value = "Psi4 (S)"
$ws.Range("A12").Font.Bold = $true

# --- Row 13: new O2 data row under the "Psi4 (S)" section ---
$ws.Range("A13").Value = "O2"
$ws.Range("B13").Value = -149.604321051363

# --- Final selection left on B16, matching the saved view state ---
$ws.Range("B16").Select() | Out-Null
